$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns F, G, H with header style copied from existing headers
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style from an existing header cell (E1) to the new header cells
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# Boolean values for rows 2-18: all FALSE except F7 which is TRUE
$boolValues = @(
    @($false, $false, $false),  # row 2
    @($false, $false, $false),  # row 3
    @($false, $false, $false),  # row 4
    @($false, $false, $false),  # row 5
    @($false, $false, $false),  # row 6
    @($true,  $false, $false),  # row 7
    @($false, $false, $false),  # row 8
    @($false, $false, $false),  # row 9
    @($false, $false, $false),  # row 10
    @($false, $false, $false),  # row 11
    @($false, $false, $false),  # row 12
    @($false, $false, $false),  # row 13
    @($false, $false, $false),  # row 14
    @($false, $false, $false),  # row 15
    @($false, $false, $false),  # row 16
    @($false, $false, $false),  # row 17
    @($false, $false, $false)   # row 18
)

for ($i = 0; $i -lt $boolValues.Length; $i++) {
    $row = $i + 2
    $vals = $boolValues[$i]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}
